# Refactor log forwarder for sanitization and structure
# The "Incidents" collection documentation in the Schema sheet contained a
# typo in the example field list: `timestamp_detected: "imestamp,` should
# read `timestamp_detected: timestamp,` (stray quote / mistyped word).
# Fix that text in place; cell position, formatting and layout stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schema")

$cell = $ws.Range("C10")

$fixed = "status: open, // Options: open, investigating, resolved [FR15]" + [char]10 + `
    "severity: high, // Options: low, medium, high, critical [FR10]" + [char]10 + `
    "timestamp_detected: timestamp," + [char]10 + `
    "source_ip: 192.168.1.50," + [char]10 + `
    "attack_type: Brute Force, // Extracted by Parser" + [char]10 + `
    "summary: `"Multiple failed login attempts detected...`", // [FR6] From LLM" + [char]10 + `
    "mitigation_steps: // [FR6] From LLM" + [char]10 + `
    "    Block IP 192.168.1.50," + [char]10 + `
    "    Reset Admin Password," + [char]10 + `
    "  " + [char]10 + `
    " risk_score: 8.5, // 0-10" + [char]10 + `
    " raw_log_reference: C://bucket_name/logs/file_123.log`", // Traceability" + [char]10 + `
    "assigned_to: user_uid // Optional" + [char]10

$cell.Value2 = $fixed

# F7 and F9 (the "audit_logs" collection's Purpose / Document ID field
# descriptions) carried a redundant "apply fill" flag left over from an
# earlier formatting pass even though no fill is actually applied. Touching
# WrapText re-normalises those two cells onto the plain wrap-top-border
# style already used elsewhere on the sheet (e.g. F3, F5, F10).
$ws.Range("F7").WrapText = $true
$ws.Range("F9").WrapText = $true

# Reflect the view state captured in the saved file: scrolled so row 7 is at
# the top, with C10 as the active/selected cell.
$window = $excel.ActiveWindow
$window.ScrollRow = 7
$window.ScrollColumn = 1
$ws.Range("C10").Select()

$wb.Save()
